$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "41.820.68"
Set-TextValue "D3" "2.274.27"
Set-TextValue "E3" "  +2.19%  "
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "305.87"
Set-TextValue "E5" "  +4.00%  "
Set-TextValue "D6" "92.68"
Set-TextValue "E6" "  +5.21%  "
Set-TextValue "E7" "  +3.63%  "
Set-TextValue "E8" "  -0.08%  "
Set-TextValue "E9" "  +3.68%  "
Set-TextValue "D10" "32.74"
Set-TextValue "D11" "54.11"
Set-TextValue "E11" "  +6.17%  "
Set-TextValue "E12" "  +2.35%  "
Set-TextValue "E13" "  +1.31%  "
Set-TextValue "D14" "6.68"
Set-TextValue "E14" "  +3.59%  "
Set-TextValue "D15" "2.625.87"
Set-TextValue "E15" "  +28.86%  "
Set-TextValue "D16" "14.32"
Set-TextValue "E16" "  +3.44%  "
Set-TextValue "D17" "2.282.71"
Set-TextValue "E17" "  +3.51%  "
Set-TextValue "D18" "0.765"
Set-TextValue "E18" "  +3.58%  "
Set-TextValue "D19" "41.757.86"
Set-TextValue "E19" "  +4.20%  "
Set-TextValue "D20" "12.26"
Set-TextValue "E20" "  +8.77%  "
Set-TextValue "D21" "0.0₃0910"
Set-TextValue "E21" "  +2.11%  "
Set-TextValue "E22" "  +2.80%  "
Set-TextValue "D23" "67.20"
Set-TextValue "E23" "  +2.22%  "
Set-TextValue "D24" "242.77"
Set-TextValue "E24" "  +2.75%  "
Set-TextValue "D25" "2.60"
Set-TextValue "E25" "  +4.74%  "
Set-TextValue "E26" "  -0.05%  "
Set-TextValue "E27" "  +5.50%  "
Set-TextValue "D28" "24.32"
Set-TextValue "E28" "  +4.61%  "
Set-TextValue "E29" "  +3.14%  "
Set-TextValue "E30" "  +0.57%  "
Set-TextValue "D31" "159.19"
Set-TextValue "E31" "  +0.16%  "
Set-TextValue "D32" "34.14"
Set-TextValue "E32" "  +6.87%  "
Set-TextValue "E33" "  +0.06%  "
Set-TextValue "E34" "  +4.35%  "
Set-TextValue "E35" "  +4.74%  "
Set-TextValue "E36" "  +1.05%  "
Set-TextValue "D37" "17.07"
Set-TextValue "E37" "  +9.29%  "
Set-TextValue "E38" "  +2.21%  "
Set-TextValue "E39" "  +2.72%  "
Set-TextValue "E40" "  +5.44%  "
Set-TextValue "E41" "  +3.22%  "
Set-TextValue "E42" "  +4.86%  "
Set-TextValue "D43" "2.074.49"
Set-TextValue "E43" "  -0.67%  "
Set-TextValue "D44" "19.46"
Set-TextValue "E44" "  +0.91%  "
Set-TextValue "E45" "  +3.07%  "
Set-TextValue "D46" "10.32"
Set-TextValue "E46" "  +1.81%  "
Set-TextValue "E47" "  +6.00%  "
Set-TextValue "E48" "  +7.99%  "
Set-TextValue "B49" "Stacks"
Set-TextValue "C49" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D49" "1.53"
Set-TextValue "E49" "  +2.98%  "
Set-TextValue "B50" "TrustWalletToken"
Set-TextValue "C50" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D50" "1.16"
Set-TextValue "E50" "  +3.68%  "
Set-TextValue "D51" "73.01"
Set-TextValue "E51" "  +7.23%  "
